$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2072
$ws.Range("J40").Value = 2176.8845
$ws.Range("L40").Value = 2176.8845
$ws.Range("N40").Value = -2526.8845

# Row 98
$ws.Range("H98").Value = 1231.25
$ws.Range("I98").Value = 1274.7916
$ws.Range("J98").Value = 1100.625
$ws.Range("K98").Value = 1274.7916
$ws.Range("L98").Value = 1100.625
$ws.Range("M98").Value = 223.2084
$ws.Range("N98").Value = -4096.625

# Row 122
$ws.Range("H122").Value = 1231.25
$ws.Range("I122").Value = 1274.7916
$ws.Range("J122").Value = 1100.625
$ws.Range("K122").Value = 3824.3748
$ws.Range("L122").Value = 3301.875
$ws.Range("M122").Value = -1374.3748
$ws.Range("N122").Value = -8201.875

# Row 129
$ws.Range("H129").Value = 1059122.1
$ws.Range("J129").Value = 1611571.5
$ws.Range("L129").Value = 4834714.5
$ws.Range("N129").Value = -4844714.5

# Row 132
$ws.Range("H132").Value = 2403.3242
$ws.Range("I132").Value = 2297.8
$ws.Range("K132").Value = 6893.400000000001
$ws.Range("M132").Value = -4363.400000000001

# Row 137
$ws.Range("H137").Value = 889
$ws.Range("I137").Value = 876.6
$ws.Range("J137").Value = 908.0769
$ws.Range("K137").Value = 2629.8
$ws.Range("L137").Value = 2724.2307
$ws.Range("M137").Value = -79.80000000000018
$ws.Range("N137").Value = -7824.2307

# Row 138
$ws.Range("H138").Value = 2576.01
$ws.Range("I138").Value = 982
$ws.Range("J138").Value = 2999.7341
$ws.Range("K138").Value = 2946
$ws.Range("L138").Value = 8999.202300000001
$ws.Range("M138").Value = 2194
$ws.Range("N138").Value = -19279.2023

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4911.5293
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# Row 74
$ws.Range("H74").Value = 912.2632
$ws.Range("I74").Value = 952.1111
$ws.Range("J74").Value = 195
$ws.Range("K74").Value = 952.1111
$ws.Range("L74").Value = 195
$ws.Range("M74").Value = -78.11109999999996
$ws.Range("N74").Value = -1943

# Row 77
$ws.Range("H77").Value = 912.2632
$ws.Range("I77").Value = 952.1111
$ws.Range("J77").Value = 195
$ws.Range("K77").Value = 4760.555499999999
$ws.Range("L77").Value = 975
$ws.Range("M77").Value = -392.5554999999995
$ws.Range("N77").Value = -9711

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 18882226
$ws.Range("I20").Value = 20837.9
$ws.Range("J20").Value = 43484036
$ws.Range("K20").Value = 20837.9
$ws.Range("L20").Value = 43484036
$ws.Range("M20").Value = -20590.9
$ws.Range("N20").Value = -43484530

# Row 94
$ws.Range("H94").Value = 1099.8572
$ws.Range("I94").Value = 1024.125
$ws.Range("J94").Value = 1342.2
$ws.Range("K94").Value = 1024.125
$ws.Range("L94").Value = 1342.2
$ws.Range("M94").Value = -573.125
$ws.Range("N94").Value = -2244.2

# Row 96
$ws.Range("H96").Value = 2976
$ws.Range("I96").Value = 2976
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2976
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0
$ws.Range("M96").Value = -230

# Row 107
$ws.Range("H107").Value = 35550
$ws.Range("I107").Value = 3133.3333
$ws.Range("K107").Value = 3133.3333
$ws.Range("M107").Value = -1213.3333

# Row 112
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954

# Row 134
$ws.Range("H134").Value = 22961.105
$ws.Range("I134").Value = 1786.1
$ws.Range("J134").Value = 60328.766
$ws.Range("K134").Value = 5358.299999999999
$ws.Range("L134").Value = 180986.298
$ws.Range("M134").Value = -2823.299999999999
$ws.Range("N134").Value = -186056.298

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3181.58
$ws.Range("I31").Value = 3233.487
$ws.Range("J31").Value = 2997.5454
$ws.Range("K31").Value = 3233.487
$ws.Range("L31").Value = 2997.5454
$ws.Range("M31").Value = -2938.487
$ws.Range("N31").Value = -3587.5454

# Row 34
$ws.Range("H34").Value = 3181.58
$ws.Range("I34").Value = 3233.487
$ws.Range("J34").Value = 2997.5454
$ws.Range("K34").Value = 3233.487
$ws.Range("L34").Value = 2997.5454
$ws.Range("M34").Value = -3031.487
$ws.Range("N34").Value = -3401.5454

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2149.7273
$ws.Range("I5").Value = 1832.9333
$ws.Range("J5").Value = 2828.5715
$ws.Range("K5").Value = 5498.7999
$ws.Range("L5").Value = 8485.7145
$ws.Range("M5").Value = -5386.7999
$ws.Range("N5").Value = -8709.7145

# Row 20
$ws.Range("H20").Value = 3800
$ws.Range("J20").Value = 3800
$ws.Range("L20").Value = 11400
$ws.Range("N20").Value = -11854

# Row 21
$ws.Range("H21").Value = 800
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 800
$ws.Range("K21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("M21").Value = 2400
$ws.Range("N21").Value = -2746

# Row 80
$ws.Range("H80").Value = 5574.1816
$ws.Range("J80").Value = 9790
$ws.Range("L80").Value = 29370
$ws.Range("N80").Value = -31242

# Row 83
$ws.Range("H83").Value = 5574.1816
$ws.Range("J83").Value = 9790
$ws.Range("L83").Value = 88110
$ws.Range("N83").Value = -97470

# Row 117
$ws.Range("H117").Value = 2260.389
$ws.Range("I117").Value = 489
$ws.Range("J117").Value = 2766.5
$ws.Range("K117").Value = 1467
$ws.Range("L117").Value = 8299.5
$ws.Range("M117").Value = 1975
$ws.Range("N117").Value = -15183.5

# Row 129
$ws.Range("H129").Value = 2463.5386
$ws.Range("I129").Value = 1280
$ws.Range("J129").Value = 3844.3333
$ws.Range("K129").Value = 3840
$ws.Range("L129").Value = 11532.9999
$ws.Range("M129").Value = 1160
$ws.Range("N129").Value = -21532.9999

# Row 135
$ws.Range("H135").Value = 2149.7273
$ws.Range("I135").Value = 1832.9333
$ws.Range("J135").Value = 2828.5715
$ws.Range("K135").Value = 16496.3997
$ws.Range("L135").Value = 25457.1435
$ws.Range("M135").Value = -13961.3997
$ws.Range("N135").Value = -30527.1435

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 12499.25
$ws.Range("J46").Value = 19000
$ws.Range("L46").Value = 19000
$ws.Range("N46").Value = -19312

# Row 70
$ws.Range("H70").Value = 4998.636
$ws.Range("J70").Value = 4998.3335
$ws.Range("L70").Value = 4998.3335
$ws.Range("N70").Value = -5538.3335

# Row 73
$ws.Range("H73").Value = 4998.636
$ws.Range("J73").Value = 4998.3335
$ws.Range("L73").Value = 4998.3335
$ws.Range("N73").Value = -6870.3335

# Row 102
$ws.Range("H102").Value = 2631.1
$ws.Range("I102").Value = 3482.4
$ws.Range("J102").Value = 1779.8
$ws.Range("K102").Value = 3482.4
$ws.Range("L102").Value = 1779.8
$ws.Range("M102").Value = -1860.4
$ws.Range("N102").Value = -5023.8

# Row 130
$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040

# Row 132
$ws.Range("H132").Value = 1996.5714
$ws.Range("I132").Value = 1782.7333
$ws.Range("J132").Value = 3279.6
$ws.Range("K132").Value = 5348.199900000001
$ws.Range("L132").Value = 9838.799999999999
$ws.Range("M132").Value = -2818.199900000001
$ws.Range("N132").Value = -14898.8

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 18705870
$ws.Range("I100").Value = 37409740
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 37409740
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -37409199
$ws.Range("N100").Value = -3082

# Row 136
$ws.Range("H136").Value = 3318.25
$ws.Range("I136").Value = 1503.5294
$ws.Range("K136").Value = 4510.5882
$ws.Range("M136").Value = -1960.5882

$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Range("H12").Value = 70007
$ws.Range("J12").Value = 70007
$ws.Range("L12").Value = 70007
$ws.Range("N12").Value = -70291

# Row 126
$ws.Range("H126").Value = 368.41666
$ws.Range("I126").Value = 244.29411
$ws.Range("J126").Value = 669.8570999999999
$ws.Range("K126").Value = 732.8823299999999
$ws.Range("L126").Value = 2009.5713
$ws.Range("M126").Value = 1737.11767
$ws.Range("N126").Value = -6949.5713
